$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where the match result ("resultado") is "Fallo" and profit is -1
$rows = @(25, 32, 33, 50, 53, 59, 61, 62, 64, 72)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "Fallo"
    $ws.Cells.Item($r, 8).Value = -1
}
